$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 275, pushing existing rows 275-288 down to 277-290
$ws.Rows.Item(275).Insert()
$ws.Rows.Item(275).Insert()

# New row 275: Espinaca "Primera" entry dated 44516
$ws.Cells.Item(275, 1).Value = 9
$ws.Cells.Item(275, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(275, 3).Value = "Metropolitana"
$ws.Cells.Item(275, 4).Value = 44516
$ws.Cells.Item(275, 5).Value = 13
$ws.Cells.Item(275, 6).Value = 100112012
$ws.Cells.Item(275, 7).Value = "Espinaca"
$ws.Cells.Item(275, 8).Value = "Sin especificar"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 250
$ws.Cells.Item(275, 11).Value = 6000
$ws.Cells.Item(275, 12).Value = 7000
$ws.Cells.Item(275, 13).Value = 6500
$ws.Cells.Item(275, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(275, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(275, 16).Value = 650
$ws.Cells.Item(275, 17).Value = 10
$ws.Cells.Item(275, 18).Value = "Hortaliza"

# New row 276: Espinaca "Segunda" entry dated 44516
$ws.Cells.Item(276, 1).Value = 9
$ws.Cells.Item(276, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(276, 3).Value = "Metropolitana"
$ws.Cells.Item(276, 4).Value = 44516
$ws.Cells.Item(276, 5).Value = 13
$ws.Cells.Item(276, 6).Value = 100112012
$ws.Cells.Item(276, 7).Value = "Espinaca"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Segunda"
$ws.Cells.Item(276, 10).Value = 97
$ws.Cells.Item(276, 11).Value = 5000
$ws.Cells.Item(276, 12).Value = 5000
$ws.Cells.Item(276, 13).Value = 5000
$ws.Cells.Item(276, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(276, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(276, 16).Value = 500
$ws.Cells.Item(276, 17).Value = 10
$ws.Cells.Item(276, 18).Value = "Hortaliza"
